# EPBDS-11327 "class" information can't be found in the OpenL datatypes
#
# Adds 4 extra test cases (case21..case24) to the existing SRClassError
# spreadsheet test-table, plus a small "Datatype MyType" block below it
# that defines a datatype with a single `value : String` field, used by
# the new test cases to exercise `.getClass()` / `.value` / `.getValue()`
# field-access resolution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- normalize the pre-existing header/label cell styles -----------------
# B4, C4, B5, C5 and B7 carried a redundant duplicate style record; re-apply
# the (identical-looking) bordered style already used by B6 so the
# workbook's style table collapses back down to its minimal form.
$ws.Range("B6").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)

# --- extend the existing table (rows 8-11) ------------------------------
# Copy the formatting of the last "real" data row (row 7) down onto the
# new rows so the new cells pick up the same borders / quote-prefix
# styling as the existing ones.
$ws.Range("B7:C7").Copy()
$ws.Range("B8:C11").PasteSpecial(-4122)

# Fill the new "Step" labels first (left column, top to bottom) ...
$ws.Range("B8").Value = "case21"
$ws.Range("B9").Value = "case22"
$ws.Range("B10").Value = "case23"
$ws.Range("B11").Value = "case24"

# --- new Datatype block (rows 14-15) -------------------------------------
$ws.Range("B14").Value = "Datatype MyType"

$ws.Range("B15").Value = "String"
$ws.Range("C15").Value = "value"
$ws.Range("D15").Value = "XXX"

# ... then go back and fill in the "Value" formulas (right column)
$ws.Range("C8").Value = "'= AccessBean.getClass()"
$ws.Range("C9").Value = "'= MyType.getClass()"
$ws.Range("C10").Value = "'= MyType.value"
$ws.Range("C11").Value = "'= MyType.getValue()"

$ws.Range("C15").Select() | Out-Null
